$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Rename headers: speedx -> speed, speedy -> score
$ws.Range("C1").Value = "speed"
$ws.Range("D1").Value = "score"

# Add new header "type" in E1
$ws.Range("E1").Value = "type"

# Move selection to F1 (was F11)
$ws.Range("F1").Select()
